$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TruckID / AssignedDockPosition / start_loading_time / end_loading_time
# rows 3-11 (row 2 and header unchanged)

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 10

$ws.Range("A4").Value = 9
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = 15

$ws.Range("A5").Value = 10
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 20

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 6

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 10

$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 5

$ws.Range("A11").Value = 6
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 10
